# Add the solution notes for the "Longest repeating character replacement"
# row: a new note in column E, row 15, matching the new shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "1. current window size - this is typically represented by a formula: (end - start + 1)"

# Leave the selection where the author ended up after entering the note.
$ws.Range("E14").Select() | Out-Null
